$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("failedTests_xe_Excel")

# Fix typo in the note for row 9 ("donet" -> "dont").
$ws.Range("O9").Value = "Initial conc. seem ok. S2 dont evolve as model"

# Clear the content (values) of rows 11-13 columns A:J, keeping formatting/styles.
$ws.Range("A11:J13").ClearContents()

# Also clear the Note column (O) for rows 12 and 13.
$ws.Range("O12:O13").ClearContents()

# Update the selection to match the committed state: entire row 11 selected,
# active cell A11.
$ws.Range("A11:XFD11").Select()
